$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at row 31 (shifts old rows 31-48 down to 32-49) ---
$ws.Rows.Item(31).Insert()

# Set values for new row 31 (new Femacal de La Calera Haba record, 2021-09-27)
$ws.Cells.Item(31,1).Value = 3
$ws.Cells.Item(31,2).Value = "Femacal de La Calera"
$ws.Cells.Item(31,3).Value = "Coquimbo"
$ws.Cells.Item(31,4).Value = 44435
$ws.Cells.Item(31,5).Value = 5
$ws.Cells.Item(31,6).Value = 100112026
$ws.Cells.Item(31,7).Value = "Haba"
$ws.Cells.Item(31,8).Value = "Sin especificar"
$ws.Cells.Item(31,9).Value = "Primera"
$ws.Cells.Item(31,10).Value = 73
$ws.Cells.Item(31,11).Value = 14000
$ws.Cells.Item(31,12).Value = 15000
$ws.Cells.Item(31,13).Value = 14521
$ws.Cells.Item(31,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(31,15).Value = "Provincia de Limarí"
$ws.Cells.Item(31,16).Value = 581
$ws.Cells.Item(31,17).Value = 25
$ws.Cells.Item(31,18).Value = "Hortaliza"

# --- Insert second new row at row 42 (shifts current rows 42-49 down to 43-50) ---
$ws.Rows.Item(42).Insert()

# Set values for new row 42 (new Femacal de La Calera Haba record, 2021-09-25)
$ws.Cells.Item(42,1).Value = 3
$ws.Cells.Item(42,2).Value = "Femacal de La Calera"
$ws.Cells.Item(42,3).Value = "Coquimbo"
$ws.Cells.Item(42,4).Value = 44433
$ws.Cells.Item(42,5).Value = 5
$ws.Cells.Item(42,6).Value = 100112026
$ws.Cells.Item(42,7).Value = "Haba"
$ws.Cells.Item(42,8).Value = "Sin especificar"
$ws.Cells.Item(42,9).Value = "Primera"
$ws.Cells.Item(42,10).Value = 73
$ws.Cells.Item(42,11).Value = 14000
$ws.Cells.Item(42,12).Value = 15000
$ws.Cells.Item(42,13).Value = 14521
$ws.Cells.Item(42,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(42,15).Value = "Provincia de Limarí"
$ws.Cells.Item(42,16).Value = 581
$ws.Cells.Item(42,17).Value = 25
$ws.Cells.Item(42,18).Value = "Hortaliza"
